$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Markov DAG state-transition pairs (column A = from-state, column B = to-state)
# matching the updated clinical-states model (Sinanovic, et al. study).
$pairs = @(
    @("Well", "Well"),
    @("Well", "Death"),
    @("Well", "Infection"),
    @("Infection", "Death"),
    @("Infection", "Infection"),
    @("Infection", "LSIL"),
    @("Infection", "HSIL"),
    @("LSIL", "Infection"),
    @("LSIL", "LSIL"),
    @("LSIL", "HSIL"),
    @("HSIL", "LSIL"),
    @("HSIL", "HSIL"),
    @("HSIL", "Infection"),
    @("HSIL", "Stage I Cancer"),
    @("Stage I Cancer", "Stage II Cancer"),
    @("Stage II Cancer", "Stage III Cancer"),
    @("Stage III Cancer", "Stage IV Cancer"),
    @("Stage I Cancer", "Treatment"),
    @("Stage II Cancer", "Treatment"),
    @("Stage III Cancer", "Treatment"),
    @("Stage IV Cancer", "Treatment"),
    @("Stage I Cancer", "Death"),
    @("Stage II Cancer", "Death"),
    @("Stage III Cancer", "Death"),
    @("Stage IV Cancer", "Death"),
    @("Treatment", "Year 1"),
    @("Year 1", "Year 2"),
    @("Year 2", "Year 3"),
    @("Year 3", "Year 4"),
    @("Year 4", "Cleared"),
    @("Cleared", "Cleared"),
    @("Cleared", "Death"),
    @("Death", "Death")
)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $pairs[$i][0]
    $ws.Cells.Item($r, 2).Value = $pairs[$i][1]
}

# Column widths (character units; engine stores width = ColumnWidth + 5/6)
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Columns.Item(2).ColumnWidth = 18.5

# Update the active selection to match the saved workbook state
$ws.Range("D28").Select()
